$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values of E2 and F2: E2 becomes 0.6, F2 becomes 0.5
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.5

# Update selection to F3
$ws.Range("F3").Select()
